$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos list with the latest scraped prices / 1h volume deltas.
# Price cells that are plain decimal numbers (e.g. "215.57") are written via
# .Formula with a leading apostrophe so Excel stores them as text (matching
# the source data, which keeps things like trailing zeros / precision that a
# numeric cell would silently drop or round).

$ws.Range("D2").Value = '25.820.44'
$ws.Range("E2").Value = '  -0.91%  '
$ws.Range("D3").Value = '1.630.01'
$ws.Range("E3").Value = '  -0.94%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Formula = '''215.57'
$ws.Range("E5").Value = '  +0.31%  '
$ws.Range("D6").Formula = '''0.5109'
$ws.Range("E6").Value = '  +0.25%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Formula = '''0.2566'
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("D9").Formula = '''0.06328'
$ws.Range("E9").Value = '  -0.44%  '
$ws.Range("D10").Formula = '''19.43'
$ws.Range("E10").Value = '  -0.72%  '
$ws.Range("D11").Formula = '''0.07777'
$ws.Range("E11").Value = '  +0.49%  '
$ws.Range("D12").Formula = '''4.238'
$ws.Range("E12").Value = '  -1.30%  '
$ws.Range("D13").Value = '1.636.05'
$ws.Range("E13").Value = '  -0.44%  '
$ws.Range("E14").Value = '  -1.07%  '
$ws.Range("D15").Formula = '''0.5511'
$ws.Range("E15").Value = '  +1.30%  '
$ws.Range("D16").Formula = '''63.61'
$ws.Range("E16").Value = '  -1.08%  '
$ws.Range("D17").Value = '0.0₅7559'
$ws.Range("E17").Value = '  -2.06%  '
$ws.Range("D18").Value = '25.855.05'
$ws.Range("E18").Value = '  -0.80%  '
$ws.Range("D19").Formula = '''1.001'
$ws.Range("E19").Value = '  +0.02%  '
$ws.Range("D20").Formula = '''194.32'
$ws.Range("E20").Value = '  -2.21%  '
$ws.Range("D21").Formula = '''4.409'
$ws.Range("E21").Value = '  -0.55%  '
$ws.Range("D22").Formula = '''9.838'
$ws.Range("E22").Value = '  -0.91%  '
$ws.Range("D23").Formula = '''6.003'
$ws.Range("E23").Value = '  -0.82%  '
$ws.Range("D24").Formula = '''1.003'
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("D25").Formula = '''1.887'
$ws.Range("E25").Value = '  +0.86%  '
$ws.Range("D26").Formula = '''142.09'
$ws.Range("E26").Value = '  +0.63%  '
$ws.Range("D27").Formula = '''0.1253'
$ws.Range("E27").Value = '  +5.02%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Formula = '''15.56'
$ws.Range("E28").Value = '  -0.34%  '
$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").Formula = '''6.746'
$ws.Range("E29").Value = '  -1.04%  '
$ws.Range("D30").Formula = '''1.237'
$ws.Range("E30").Value = '  +0.18%  '
$ws.Range("D31").Formula = '''0.04885'
$ws.Range("E31").Value = '  +0.61%  '
$ws.Range("D32").Formula = '''3.233'
$ws.Range("E32").Value = '  -0.82%  '
$ws.Range("D33").Formula = '''3.173'
$ws.Range("E33").Value = '  +0.18%  '
$ws.Range("D34").Formula = '''1.544'
$ws.Range("E34").Value = '  +1.08%  '
$ws.Range("D35").Formula = '''2.373'
$ws.Range("E35").Value = '  +0.30%  '
$ws.Range("D36").Formula = '''0.8940'
$ws.Range("E36").Value = '  -0.65%  '
$ws.Range("D37").Formula = '''0.5516'
$ws.Range("E37").Value = '  +0.92%  '
$ws.Range("D38").Formula = '''2.539'
$ws.Range("E38").Value = '  -1.77%  '
$ws.Range("D39").Value = '1.113.86'
$ws.Range("E39").Value = '  -2.59%  '
$ws.Range("D40").Formula = '''0.01550'
$ws.Range("E40").Value = '  -0.92%  '
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("D42").Formula = '''5.569'
$ws.Range("E42").Value = '  +3.35%  '
$ws.Range("D43").Formula = '''0.7951'
$ws.Range("E43").Value = '  -2.13%  '
$ws.Range("D44").Formula = '''97.39'
$ws.Range("E44").Value = '  -1.96%  '
$ws.Range("D45").Value = '1.777.46'
$ws.Range("E45").Value = '  -0.28%  '
$ws.Range("E46").Value = '  -10.41%  '
$ws.Range("D47").Formula = '''0.4436'
$ws.Range("E47").Value = '  -2.07%  '
$ws.Range("D48").Formula = '''1.001'
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("D49").Formula = '''54.60'
$ws.Range("E49").Value = '  -0.67%  '
$ws.Range("D50").Formula = '''0.05134'
$ws.Range("E50").Value = '  +1.52%  '
$ws.Range("D51").Formula = '''7.555'
$ws.Range("E51").Value = '  +3.21%  '
